$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where Total Attendance Count (D) and Real (E) become 1
$dePresentRows = @(4, 5, 9, 10, 11)
foreach ($r in $dePresentRows) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
}

# Rows where Absent (H) becomes 1
$absentRows = @(3, 6, 7, 8, 12, 13, 14, 15, 16, 17, 18)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1
}
